# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Updates the "Estado de Cuenta" sheet:
#   - Bumps the total "VALOR MORA" (E11) and "Cant. Periodos" (F13)
#   - Swaps the "Novedad de Ingreso" / "Novedad de Retiro" column headers (H15/I15)
#   - Adds a new detail row (18) for period 2509, cloned from the existing
#     period-2508 row (17), pushing the signature block down one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the total overdue amount shown at the top of the statement.
$ws.Range("E11").Value = 145880

# 2. One more period is now being reported.
$ws.Range("F13").Value = 3

# 3. The "Novedad de Ingreso" / "Novedad de Retiro" headers were swapped.
$novIngreso = $ws.Range("H15").Value()
$novRetiro = $ws.Range("I15").Value()
$ws.Range("H15").Value = $novRetiro
$ws.Range("I15").Value = $novIngreso

# 4. Insert a new detail row right after the existing 2508 row, cloning its
#    formatting, then point it at the new 2509 period with the same worker.
$ws.Rows.Item(18).Insert()

$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

$ws.Range("B18").Value = $ws.Range("B17").Value()
$ws.Range("C18").Value = $ws.Range("C17").Value()
$ws.Range("D18").Value = $ws.Range("D17").Value()
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500
